$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update quarter header labels (columns E..N) on all 6 header rows ---
$newHeaders = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)
$headerCols = @("E","F","G","H","I","J","K","L","M","N")
$headerRows = @(8,21,35,46,56,67)
foreach ($r in $headerRows) {
    for ($i = 0; $i -lt 10; $i++) {
        $ws.Range($headerCols[$i] + $r).Value = $newHeaders[$i]
    }
}

# --- Update data rows: shift quarters left (E<-F_old..M<-N_old) and set new N with latest quarter data ---
# Row 10
$ws.Range("E10").Value = "-"
$ws.Range("F10").Value = "-"
$ws.Range("G10").Value = "-"
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0

# Row 11
$ws.Range("E11").Value = 4799418
$ws.Range("F11").Value = 4561797
$ws.Range("G11").Value = 3627810
$ws.Range("H11").Value = "-"
$ws.Range("I11").Value = "-"
$ws.Range("J11").Value = 2509811
$ws.Range("K11").Value = 5275511
$ws.Range("L11").Value = 4674724
$ws.Range("M11").Value = 4357531
$ws.Range("N11").Value = -5459031

# Row 12
$ws.Range("E12").Value = "-"
$ws.Range("F12").Value = "-"
$ws.Range("G12").Value = "-"
$ws.Range("H12").Value = 2727610
$ws.Range("I12").Value = 3778799
$ws.Range("J12").Value = "-"
$ws.Range("K12").Value = "-"
$ws.Range("L12").Value = "-"
$ws.Range("M12").Value = "-"
$ws.Range("N12").Value = "-"

# Row 13
$ws.Range("E13").Value = 10991416
$ws.Range("F13").Value = 10179392
$ws.Range("G13").Value = 6339988
$ws.Range("H13").Value = "-"
$ws.Range("I13").Value = "-"
$ws.Range("J13").Value = 8543605
$ws.Range("K13").Value = 9698138
$ws.Range("L13").Value = 10735735
$ws.Range("M13").Value = 10103272
$ws.Range("N13").Value = 10485836

# Row 14
$ws.Range("E14").Value = "-"
$ws.Range("F14").Value = "-"
$ws.Range("G14").Value = "-"
$ws.Range("H14").Value = 9767264
$ws.Range("I14").Value = 8283461
$ws.Range("J14").Value = "-"
$ws.Range("K14").Value = "-"
$ws.Range("L14").Value = "-"
$ws.Range("M14").Value = "-"
$ws.Range("N14").Value = "-"

# Row 15
$ws.Range("E15").Value = 2764458
$ws.Range("F15").Value = 2296270
$ws.Range("G15").Value = 1478935
$ws.Range("H15").Value = "-"
$ws.Range("I15").Value = "-"
$ws.Range("J15").Value = 2473137
$ws.Range("K15").Value = 1497447
$ws.Range("L15").Value = 2506105
$ws.Range("M15").Value = 2534778
$ws.Range("N15").Value = 12329583

# Row 16
$ws.Range("E16").Value = "-"
$ws.Range("F16").Value = "-"
$ws.Range("G16").Value = "-"
$ws.Range("H16").Value = 2043778
$ws.Range("I16").Value = 1935946
$ws.Range("J16").Value = "-"
$ws.Range("K16").Value = "-"
$ws.Range("L16").Value = "-"
$ws.Range("M16").Value = "-"
$ws.Range("N16").Value = "-"

# Row 17
$ws.Range("E17").Value = 18555292
$ws.Range("F17").Value = 17037459
$ws.Range("G17").Value = 11446733
$ws.Range("H17").Value = 14538652
$ws.Range("I17").Value = 13998206
$ws.Range("J17").Value = 13526553
$ws.Range("K17").Value = 16471096
$ws.Range("L17").Value = 17916564
$ws.Range("M17").Value = 16995581
$ws.Range("N17").Value = 17356388

# Row 23
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = "-"
$ws.Range("G23").Value = "-"
$ws.Range("H23").Value = "-"
$ws.Range("I23").Value = "-"
$ws.Range("J23").Value = "-"
$ws.Range("K23").Value = "-"
$ws.Range("L23").Value = "-"
$ws.Range("M23").Value = "-"
$ws.Range("N23").Value = "-"

# Row 24
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = "-"
$ws.Range("G24").Value = "-"
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0
$ws.Range("N24").Value = 0

# Row 25
$ws.Range("E25").Value = 3574712
$ws.Range("F25").Value = 3902178
$ws.Range("G25").Value = 3503475
$ws.Range("H25").Value = "-"
$ws.Range("I25").Value = "-"
$ws.Range("J25").Value = 3526832
$ws.Range("K25").Value = 3847863
$ws.Range("L25").Value = 4509662
$ws.Range("M25").Value = 2921760
$ws.Range("N25").Value = 4485770

# Row 26
$ws.Range("E26").Value = "-"
$ws.Range("F26").Value = "-"
$ws.Range("G26").Value = "-"
$ws.Range("H26").Value = 2653328
$ws.Range("I26").Value = 3969266
$ws.Range("J26").Value = "-"
$ws.Range("K26").Value = "-"
$ws.Range("L26").Value = "-"
$ws.Range("M26").Value = "-"
$ws.Range("N26").Value = "-"

# Row 27
$ws.Range("E27").Value = 8114182
$ws.Range("F27").Value = 9657650
$ws.Range("G27").Value = 6904740
$ws.Range("H27").Value = "-"
$ws.Range("I27").Value = "-"
$ws.Range("J27").Value = 10466061
$ws.Range("K27").Value = 9747934
$ws.Range("L27").Value = 10012442
$ws.Range("M27").Value = 8626882
$ws.Range("N27").Value = 9221271

# Row 28
$ws.Range("E28").Value = "-"
$ws.Range("F28").Value = "-"
$ws.Range("G28").Value = "-"
$ws.Range("H28").Value = 11045647
$ws.Range("I28").Value = 7334750
$ws.Range("J28").Value = "-"
$ws.Range("K28").Value = "-"
$ws.Range("L28").Value = "-"
$ws.Range("M28").Value = "-"
$ws.Range("N28").Value = "-"

# Row 29
$ws.Range("E29").Value = 1576384
$ws.Range("F29").Value = 3225922
$ws.Range("G29").Value = 1973674
$ws.Range("H29").Value = "-"
$ws.Range("I29").Value = "-"
$ws.Range("J29").Value = 2300045
$ws.Range("K29").Value = 3063682
$ws.Range("L29").Value = 1584239
$ws.Range("M29").Value = 6516995
$ws.Range("N29").Value = 3226290

# Row 30
$ws.Range("E30").Value = "-"
$ws.Range("F30").Value = "-"
$ws.Range("G30").Value = "-"
$ws.Range("H30").Value = 2019843
$ws.Range("I30").Value = 2935819
$ws.Range("J30").Value = "-"
$ws.Range("K30").Value = "-"
$ws.Range("L30").Value = "-"
$ws.Range("M30").Value = "-"
$ws.Range("N30").Value = "-"

# Row 31
$ws.Range("E31").Value = 13265278
$ws.Range("F31").Value = 16785750
$ws.Range("G31").Value = 12381889
$ws.Range("H31").Value = 15718818
$ws.Range("I31").Value = 14239835
$ws.Range("J31").Value = 16292938
$ws.Range("K31").Value = 16659479
$ws.Range("L31").Value = 16106343
$ws.Range("M31").Value = 18065637
$ws.Range("N31").Value = 16933331

# Row 37
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = "-"
$ws.Range("G37").Value = "-"
$ws.Range("H37").Value = "-"
$ws.Range("I37").Value = "-"
$ws.Range("J37").Value = "-"
$ws.Range("K37").Value = "-"
$ws.Range("L37").Value = "-"
$ws.Range("M37").Value = "-"
$ws.Range("N37").Value = "-"

# Row 38
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = "-"
$ws.Range("G38").Value = "-"
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 0
$ws.Range("N38").Value = 0

# Row 39
$ws.Range("E39").Value = 824591
$ws.Range("F39").Value = 1080895
$ws.Range("G39").Value = 1060392
$ws.Range("H39").Value = 943321
$ws.Range("I39").Value = 1142457
$ws.Range("J39").Value = 1111544
$ws.Range("K39").Value = 1194903
$ws.Range("L39").Value = 2158823
$ws.Range("M39").Value = 279891
$ws.Range("N39").Value = 1440889

# Row 40
$ws.Range("E40").Value = 1721847
$ws.Range("F40").Value = 2292871
$ws.Range("G40").Value = 1581151
$ws.Range("H40").Value = 2378942
$ws.Range("I40").Value = 1882454
$ws.Range("J40").Value = 2697091
$ws.Range("K40").Value = 2733144
$ws.Range("L40").Value = 3470136
$ws.Range("M40").Value = 3120662
$ws.Range("N40").Value = 4568718

# Row 41
$ws.Range("E41").Value = 585463
$ws.Range("F41").Value = 1021372
$ws.Range("G41").Value = 748168
$ws.Range("H41").Value = 958207
$ws.Range("I41").Value = 1328298
$ws.Range("J41").Value = 1378084
$ws.Range("K41").Value = 1691608
$ws.Range("L41").Value = 1015053
$ws.Range("M41").Value = 3919820
$ws.Range("N41").Value = 1667080

# Row 42
$ws.Range("E42").Value = 3131901
$ws.Range("F42").Value = 4395138
$ws.Range("G42").Value = 3389711
$ws.Range("H42").Value = 4280470
$ws.Range("I42").Value = 4353209
$ws.Range("J42").Value = 5186719
$ws.Range("K42").Value = 5619655
$ws.Range("L42").Value = 6644012
$ws.Range("M42").Value = 7320373
$ws.Range("N42").Value = 7676687

# Row 48
$ws.Range("E48").Value = "-"
$ws.Range("F48").Value = "-"
$ws.Range("G48").Value = "-"
$ws.Range("H48").Value = "-"
$ws.Range("I48").Value = "-"
$ws.Range("J48").Value = "-"
$ws.Range("K48").Value = "-"
$ws.Range("L48").Value = "-"
$ws.Range("M48").Value = "-"
$ws.Range("N48").Value = "-"

# Row 49
$ws.Range("E49").Value = "-"
$ws.Range("F49").Value = "-"
$ws.Range("G49").Value = "-"
$ws.Range("H49").Value = "-"
$ws.Range("I49").Value = "-"
$ws.Range("J49").Value = "-"
$ws.Range("K49").Value = "-"
$ws.Range("L49").Value = "-"
$ws.Range("M49").Value = "-"
$ws.Range("N49").Value = "-"

# Row 50
$ws.Range("E50").Value = 230673
$ws.Range("F50").Value = 276998
$ws.Range("G50").Value = 302669
$ws.Range("H50").Value = 355524
$ws.Range("I50").Value = 276263
$ws.Range("J50").Value = 315168
$ws.Range("K50").Value = 310537
$ws.Range("L50").Value = 478711
$ws.Range("M50").Value = 95795
$ws.Range("N50").Value = 321213

# Row 51
$ws.Range("E51").Value = 212202
$ws.Range("F51").Value = 237415
$ws.Range("G51").Value = 228995
$ws.Range("H51").Value = 215374
$ws.Range("I51").Value = 288061
$ws.Range("J51").Value = 257699
$ws.Range("K51").Value = 280382
$ws.Range("L51").Value = 346582
$ws.Range("M51").Value = 361737
$ws.Range("N51").Value = 495454

# Row 52
$ws.Range("E52").Value = 371396
$ws.Range("F52").Value = 316614
$ws.Range("G52").Value = 379074
$ws.Range("H52").Value = 474397
$ws.Range("I52").Value = 403431
$ws.Range("J52").Value = 599155
$ws.Range("K52").Value = 552149
$ws.Range("L52").Value = 640720
$ws.Range("M52").Value = 601477
$ws.Range("N52").Value = 516717

# Row 58
$ws.Range("E58").Value = 0
$ws.Range("F58").Value = "-"
$ws.Range("G58").Value = "-"
$ws.Range("H58").Value = "-"
$ws.Range("I58").Value = "-"
$ws.Range("J58").Value = "-"
$ws.Range("K58").Value = "-"
$ws.Range("L58").Value = "-"
$ws.Range("M58").Value = "-"
$ws.Range("N58").Value = "-"

# Row 59
$ws.Range("E59").Value = 0
$ws.Range("F59").Value = "-"
$ws.Range("G59").Value = "-"
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = 0
$ws.Range("N59").Value = 0

# Row 60
$ws.Range("E60").Value = -273146
$ws.Range("F60").Value = -611085
$ws.Range("G60").Value = -392752
$ws.Range("H60").Value = -440652
$ws.Range("I60").Value = -433412
$ws.Range("J60").Value = -646733
$ws.Range("K60").Value = -575481
$ws.Range("L60").Value = -988695
$ws.Range("M60").Value = -499159
$ws.Range("N60").Value = -926933

# Row 61
$ws.Range("E61").Value = -771148
$ws.Range("F61").Value = -1533587
$ws.Range("G61").Value = -876056
$ws.Range("H61").Value = -1247519
$ws.Range("I61").Value = -1212546
$ws.Range("J61").Value = -1807179
$ws.Range("K61").Value = -1618129
$ws.Range("L61").Value = -2120600
$ws.Range("M61").Value = -1572074
$ws.Range("N61").Value = -3164387

# Row 62
$ws.Range("E62").Value = -242505
$ws.Range("F62").Value = -575404
$ws.Range("G62").Value = -416065
$ws.Range("H62").Value = -437091
$ws.Range("I62").Value = -566928
$ws.Range("J62").Value = -731857
$ws.Range("K62").Value = -883913
$ws.Range("L62").Value = -304078
$ws.Range("M62").Value = -2110977
$ws.Range("N62").Value = -399633

# Row 63
$ws.Range("E63").Value = -1286799
$ws.Range("F63").Value = -2720076
$ws.Range("G63").Value = -1684873
$ws.Range("H63").Value = -2125262
$ws.Range("I63").Value = -2212886
$ws.Range("J63").Value = -3185769
$ws.Range("K63").Value = -3077523
$ws.Range("L63").Value = -3413373
$ws.Range("M63").Value = -4182210
$ws.Range("N63").Value = -4490953

# Row 69
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = "-"
$ws.Range("G69").Value = "-"
$ws.Range("H69").Value = "-"
$ws.Range("I69").Value = "-"
$ws.Range("J69").Value = "-"
$ws.Range("K69").Value = "-"
$ws.Range("L69").Value = "-"
$ws.Range("M69").Value = "-"
$ws.Range("N69").Value = "-"

# Row 70
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = "-"
$ws.Range("G70").Value = "-"
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = 0
$ws.Range("N70").Value = 0

# Row 71
$ws.Range("E71").Value = 551445
$ws.Range("F71").Value = 469809
$ws.Range("G71").Value = 667640
$ws.Range("H71").Value = 502669
$ws.Range("I71").Value = 663151
$ws.Range("J71").Value = 464811
$ws.Range("K71").Value = 619422
$ws.Range("L71").Value = 1170128
$ws.Range("M71").Value = -219268
$ws.Range("N71").Value = 513956

# Row 72
$ws.Range("E72").Value = 950699
$ws.Range("F72").Value = 759285
$ws.Range("G72").Value = 705095
$ws.Range("H72").Value = 1131423
$ws.Range("I72").Value = 900308
$ws.Range("J72").Value = 889912
$ws.Range("K72").Value = 1115015
$ws.Range("L72").Value = 1349536
$ws.Range("M72").Value = 1548588
$ws.Range("N72").Value = 1404331

# Row 73
$ws.Range("E73").Value = 342958
$ws.Range("F73").Value = 445968
$ws.Range("G73").Value = 332103
$ws.Range("H73").Value = 521116
$ws.Range("I73").Value = 576864
$ws.Range("J73").Value = 646227
$ws.Range("K73").Value = 807695
$ws.Range("L73").Value = 710975
$ws.Range("M73").Value = 1808843
$ws.Range("N73").Value = 1267447

# Row 74
$ws.Range("E74").Value = 1845102
$ws.Range("F74").Value = 1675062
$ws.Range("G74").Value = 1704838
$ws.Range("H74").Value = 2155208
$ws.Range("I74").Value = 2140323
$ws.Range("J74").Value = 2000950
$ws.Range("K74").Value = 2542132
$ws.Range("L74").Value = 3230639
$ws.Range("M74").Value = 3138163
$ws.Range("N74").Value = 3185734
